# Update "want to go" counts (column F) for specific rows on the sheets
# that hold the full event list ("展览" and "全部类型" sheets, i.e. the
# 1st and 4th worksheet in this workbook). The other two sheets ("演出"
# and "本地生活") only contain a header row and are left untouched.

$wb = $excel.ActiveWorkbook

# Mapping of row number -> new value for column F
$updates = @{
    3  = 1318
    8  = 12
    11 = 4437
    12 = 6706
    18 = 4098
    19 = 464
    22 = 2676
    24 = 545
    26 = 343
    27 = 347
    31 = 1611
    35 = 76
    36 = 536
    37 = 494
    40 = 624
}

$sheetIndexes = @(1, 4)

foreach ($sheetIndex in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    foreach ($row in $updates.Keys) {
        $newValue = $updates[$row]
        $cellAddress = "F" + $row
        $ws.Range($cellAddress).Value = $newValue
    }
}
